$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-13 Sunday" "2025-04-14 Monday"

Replace-Text "949÷4=237, 1" "396÷4=99, 0"
Replace-Text "121÷5=24, 1" "274÷2=137, 0"
Replace-Text "564÷7=80, 4" "360÷4=90, 0"
Replace-Text "102÷2=51, 0" "921÷9=102, 3"
Replace-Text "902÷9=100, 2" "315÷7=45, 0"

Replace-Text "747÷2=373, 1" "734÷6=122, 2"
Replace-Text "570÷7=81, 3" "582÷8=72, 6"
Replace-Text "228÷2=114, 0" "515÷2=257, 1"
Replace-Text "901÷5=180, 1" "914÷3=304, 2"
Replace-Text "942÷2=471, 0" "270÷5=54, 0"

Replace-Text "311÷9=34, 5" "557÷5=111, 2"
Replace-Text "127÷4=31, 3" "755÷9=83, 8"
Replace-Text "546÷3=182, 0" "468÷2=234, 0"
Replace-Text "169÷9=18, 7" "220÷6=36, 4"
Replace-Text "487÷5=97, 2" "573÷7=81, 6"

Replace-Text "516÷8=64, 4" "993÷2=496, 1"
Replace-Text "425÷9=47, 2" "930÷6=155, 0"
Replace-Text "506÷2=253, 0" "680÷6=113, 2"
Replace-Text "299÷8=37, 3" "930÷4=232, 2"
Replace-Text "758÷4=189, 2" "478÷2=239, 0"

Replace-Text "951÷8=118, 7" "790÷5=158, 0"
Replace-Text "406÷9=45, 1" "491÷8=61, 3"
Replace-Text "343÷3=114, 1" "723÷4=180, 3"
Replace-Text "223÷4=55, 3" "380÷6=63, 2"
Replace-Text "728÷6=121, 2" "759÷8=94, 7"
